$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.905.74'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +2.21%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.813.44'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +2.96%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.006'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.44%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.74'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +3.62%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.005'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.43%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4297'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.68%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3698'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +2.31%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07241'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +2.91%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8659'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +4.33%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '2.051.66'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +19.91%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +5.91%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.642'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +4.14%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.401'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +3.14%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06936'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +2.18%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '80.88'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +2.28%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.006'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.41%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008836'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +2.19%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.58%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.20'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.94%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '26.936.36'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +4.22%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.199'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +4.20%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.97'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.08%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.300.83'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +17.97%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '154.32'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.53%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.884'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -1.12%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.32'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +1.02%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.243'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +4.65%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.923'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +15.53%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '114.62'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.09%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.09%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7442'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +3.05%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.162'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +4.54%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.431'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +3.16%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.803'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +3.20%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.53%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.123'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +5.43%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05232'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +2.91%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01927'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.64%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5101'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +4.39%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.747'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +10.50%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1653'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +3.26%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.478'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +5.36%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.295'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +4.16%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '107.43'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +2.86%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.44'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +4.58%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.005'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.47%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.651'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +5.49%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4563'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +2.45%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06270'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.43%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.803'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +5.76%  '
